# Updated cryptos list on Thu Feb 22 15:55:47 UTC 2024 with GitHub Actions
# Refresh the scraped Price (D) / Volume(1h) (E) columns with the latest
# values, and rotate the 3 rows whose ranking order changed (Hedera/Dai/
# EthereumClassic, rows 29-31) while keeping their rank index (col A) fixed.
#
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as plain text (matching the original inlineStr cells)
# instead of auto-coercing them to numbers; .Style is reset to "Normal"
# right after so the quote-prefix marker doesn't linger as a format change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.500.85"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "2.980.71"
$ws.Range("E3").Value = "  +2.41%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'379.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.97%  "

$ws.Range("D6").Value = "'105.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "

$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.594"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").Value = "'37.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "

$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "'0.0847"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").Value = "3.449.38"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("D14").Value = "'18.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "'7.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").Value = "2.977.76"
$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("E17").Value = "  +2.84%  "

$ws.Range("D18").Value = "51.486.87"
$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("D19").Value = "'3.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.13%  "

$ws.Range("D20").Value = "'7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.87%  "

$ws.Range("D21").Value = "'12.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("E22").Value = "  +1.98%  "

$ws.Range("D23").Value = "'69.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.80%  "

$ws.Range("D24").Value = "'262.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").Value = "'2.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.00%  "

$ws.Range("D26").Value = "'7.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +23.55%  "

$ws.Range("D27").Value = "'7.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.96%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'25.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.112"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.45%  "

$ws.Range("D32").Value = "'9.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("D33").Value = "'35.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "

$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("D35").Value = "'50.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").Value = "'0.0447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.14%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").Value = "'17.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").Value = "'2.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.37%  "

$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("D43").Value = "'124.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.38%  "

$ws.Range("D44").Value = "'22.09"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.291"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.55%  "

$ws.Range("D46").Value = "'2.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.42%  "

$ws.Range("E47").Value = "  +2.83%  "

$ws.Range("D48").Value = "2.046.76"
$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("D49").Value = "'3.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D51").Value = "'5.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.71%  "
